# This workbook contains a weekly price log for "Jengibre" (ginger) at the
# "Terminal La Palmera de La Serena" market. A new weekly record is inserted
# right after the existing row for 2023-09-22 (row 63), pushing every
# subsequent record down by one row and extending the sheet from 203 to 204
# data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 63; this shifts rows 63..203 down to
# 64..204 and extends the used range to A1:R204.
$ws.Rows(63).Insert()

# Populate the newly inserted row 63 with the new weekly record.
$ws.Cells.Item(63, 1).Value2  = 8
$ws.Cells.Item(63, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(63, 3).Value2  = "Coquimbo"
$ws.Cells.Item(63, 4).Value2  = 45251
$ws.Cells.Item(63, 5).Value2  = 4
$ws.Cells.Item(63, 6).Value2  = 100114007
$ws.Cells.Item(63, 7).Value2  = "Jengibre"
$ws.Cells.Item(63, 8).Value2  = "Sin especificar"
$ws.Cells.Item(63, 9).Value2  = "Primera"
$ws.Cells.Item(63, 10).Value2 = 360
$ws.Cells.Item(63, 11).Value2 = 24000
$ws.Cells.Item(63, 12).Value2 = 25000
$ws.Cells.Item(63, 13).Value2 = 24500
$ws.Cells.Item(63, 14).Value2 = "$/caja 13 kilos"
$ws.Cells.Item(63, 15).Value2 = "Perú"
$ws.Cells.Item(63, 16).Value2 = 1885
$ws.Cells.Item(63, 17).Value2 = 13
$ws.Cells.Item(63, 18).Value2 = "Hortaliza"
